# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the ee2b7eb3... file row
# on both locale sheets (zh-cn and de-de), reflecting a freshly regenerated
# handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-11 07:47:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-11 07:47:40"
